# Auto-generated edit script applying the Ridill_Profits.xlsx diff
# Updates H:N numeric columns across ALC, ARM, CRP, CUL, GSM, LTW, WVR sheets
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 5
$ws.Range("H5").Value = 68.2
$ws.Range("I5").Value = 60.25
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 60.25
$ws.Range("L5").Value = 100
$ws.Range("M5").Value = 54.75
$ws.Range("N5").Value = -330
# Row 96
$ws.Range("H96").Value = 941.6667
$ws.Range("I96").Value = 438.57144
$ws.Range("J96").Value = 1646
$ws.Range("K96").Value = 1315.71432
$ws.Range("L96").Value = 4938
$ws.Range("M96").Value = 57.28567999999996
$ws.Range("N96").Value = -7684
# Row 116
$ws.Range("H116").Value = 12340981
$ws.Range("I116").Value = 6669516.5
$ws.Range("K116").Value = 6669516.5
$ws.Range("M116").Value = -6666074.5

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 37
$ws.Range("H37").Value = 12303.077
$ws.Range("I37").Value = 6980
$ws.Range("K37").Value = 6980
$ws.Range("M37").Value = -6707

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 1401
$ws.Range("I4").Value = 1401
$ws.Range("K4").Value = 1401
$ws.Range("M4").Value = -1289
# Row 50
$ws.Range("H50").Value = 13818.429
$ws.Range("J50").Value = 13818.429
$ws.Range("L50").Value = 13818.429
$ws.Range("N50").Value = -15068.429
# Row 51
$ws.Range("H51").Value = 9329
$ws.Range("J51").Value = 9329
$ws.Range("L51").Value = 9329
$ws.Range("N51").Value = -10801
# Row 59
$ws.Range("H59").Value = 16501
$ws.Range("J59").Value = 16501
$ws.Range("L59").Value = 16501
$ws.Range("N59").Value = -18791
# Row 60
$ws.Range("H60").Value = 3800.2856
$ws.Range("J60").Value = 8301
$ws.Range("L60").Value = 8301
$ws.Range("N60").Value = -9323
# Row 61
$ws.Range("H61").Value = 9329
$ws.Range("J61").Value = 9329
$ws.Range("L61").Value = 9329
$ws.Range("N61").Value = -10025
# Row 70
$ws.Range("H70").Value = 17980
$ws.Range("J70").Value = 17980
$ws.Range("L70").Value = 17980
$ws.Range("N70").Value = -18610
# Row 73
$ws.Range("H73").Value = 17980
$ws.Range("J73").Value = 17980
$ws.Range("L73").Value = 17980
$ws.Range("N73").Value = -20164
# Row 74
$ws.Range("H74").Value = 16137.167
$ws.Range("J74").Value = 16137.167
$ws.Range("L74").Value = 16137.167
$ws.Range("N74").Value = -17885.167
# Row 77
$ws.Range("H77").Value = 16137.167
$ws.Range("J77").Value = 16137.167
$ws.Range("L77").Value = 48411.501
$ws.Range("N77").Value = -57147.501
# Row 86
$ws.Range("H86").Value = 11193.143
$ws.Range("I86").Value = 12754.909
$ws.Range("J86").Value = 5466.6665
$ws.Range("K86").Value = 12754.909
$ws.Range("L86").Value = 5466.6665
$ws.Range("M86").Value = -11631.909
$ws.Range("N86").Value = -7712.6665
# Row 89
$ws.Range("H89").Value = 11193.143
$ws.Range("I89").Value = 12754.909
$ws.Range("J89").Value = 5466.6665
$ws.Range("K89").Value = 63774.545
$ws.Range("L89").Value = 27333.3325
$ws.Range("M89").Value = -58158.545
$ws.Range("N89").Value = -38565.3325

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 48.5
$ws.Range("I2").Value = 67.5
$ws.Range("J2").Value = 29.5
$ws.Range("K2").Value = 405
$ws.Range("L2").Value = 177
$ws.Range("M2").Value = -292
$ws.Range("N2").Value = -403
# Row 4
$ws.Range("H4").Value = 387.75
$ws.Range("I4").Value = 99.59999999999999
$ws.Range("J4").Value = 868
$ws.Range("K4").Value = 298.8
$ws.Range("L4").Value = 2604
$ws.Range("M4").Value = -186.8
$ws.Range("N4").Value = -2828
# Row 5
$ws.Range("H5").Value = 9959939
$ws.Range("I5").Value = 15385596
$ws.Range("J5").Value = 6945684.5
$ws.Range("K5").Value = 46156788
$ws.Range("L5").Value = 20837053.5
$ws.Range("M5").Value = -46156676
$ws.Range("N5").Value = -20837277.5
# Row 107
$ws.Range("H107").Value = 1222075.2
$ws.Range("J107").Value = 1127.75
$ws.Range("L107").Value = 3383.25
$ws.Range("N107").Value = -7223.25
# Row 122
$ws.Range("H122").Value = 540.8
$ws.Range("I122").Value = 279.09525
$ws.Range("K122").Value = 2511.85725
$ws.Range("M122").Value = -61.85725000000002
# Row 135
$ws.Range("H135").Value = 9959939
$ws.Range("I135").Value = 15385596
$ws.Range("J135").Value = 6945684.5
$ws.Range("K135").Value = 138470364
$ws.Range("L135").Value = 62511160.5
$ws.Range("M135").Value = -138467829
$ws.Range("N135").Value = -62516230.5

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 25001100
$ws.Range("I97").Value = 1300
$ws.Range("J97").Value = 50000896
$ws.Range("K97").Value = 1300
$ws.Range("L97").Value = 50000896
$ws.Range("M97").Value = -804
$ws.Range("N97").Value = -50001888
# Row 126
$ws.Range("H126").Value = 7068.2383
$ws.Range("I126").Value = 8996.933999999999
$ws.Range("J126").Value = 2246.5
$ws.Range("K126").Value = 26990.802
$ws.Range("L126").Value = 6739.5
$ws.Range("M126").Value = -24520.802
$ws.Range("N126").Value = -11679.5

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()
# Row 7
$ws.Range("H7").Value = 1914
$ws.Range("I7").Value = 1448.5714
$ws.Range("J7").Value = 3000
$ws.Range("K7").Value = 1448.5714
$ws.Range("L7").Value = 3000
$ws.Range("M7").Value = -1336.5714
$ws.Range("N7").Value = -3224
# Row 61
$ws.Range("H61").Value = 2273.8667
$ws.Range("I61").Value = 900.8889
$ws.Range("J61").Value = 4333.3335
$ws.Range("K61").Value = 900.8889
$ws.Range("L61").Value = 4333.3335
$ws.Range("M61").Value = -698.8889
$ws.Range("N61").Value = -4737.3335
# Row 113
$ws.Range("H113").Value = 2273.8667
$ws.Range("I113").Value = 900.8889
$ws.Range("J113").Value = 4333.3335
$ws.Range("K113").Value = 900.8889
$ws.Range("L113").Value = 4333.3335
$ws.Range("M113").Value = 1269.1111
$ws.Range("N113").Value = -8673.333500000001
# Row 122
$ws.Range("H122").Value = 17023770
$ws.Range("I122").Value = 1775750
$ws.Range("K122").Value = 5327250
$ws.Range("M122").Value = -5324800
# Row 126
$ws.Range("H126").Value = 1914
$ws.Range("I126").Value = 1448.5714
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 4345.7142
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -1875.7142
$ws.Range("N126").Value = -13940

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 2209.6155
$ws.Range("I122").Value = 2068
$ws.Range("J122").Value = 2298.125
$ws.Range("K122").Value = 6204
$ws.Range("L122").Value = 6894.375
$ws.Range("M122").Value = -3754
$ws.Range("N122").Value = -11794.375
